$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the value in A1 from "Ayaansh" to "Janaki"
$ws.Range("A1").Value = "Janaki"
$ws.Range("A1").Select()
